$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
# D1 header text changes from "Oct. Deployments" to "Oct. Days Deployed"
$ws.Range("D1").Value = "Oct. Days Deployed"

# I1 / J1 headers were swapped ("MTTR Sept." <-> "Failed Changes Oct.")
$i1 = $ws.Range("I1").Value2
$j1 = $ws.Range("J1").Value2
$ws.Range("I1").Value = $j1
$ws.Range("J1").Value = $i1

# --- Updated metric values for the "Oct. Days Deployed" (D) and "Sep. Deployments" (E) columns ---
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 2

$ws.Range("E2").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 3

# --- New average formulas for columns D and E, matching existing F8/G8 formatting ---
$ws.Range("D8").Formula = "=AVERAGE(D2:D7)"
$ws.Range("E8").Formula = "=AVERAGE(E2:E7)"
$ws.Range("D8").NumberFormat = "0"
$ws.Range("E8").NumberFormat = "0"

# --- Page setup: switch to portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Move the active selection to E9 ---
$ws.Range("E9").Select()
